$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (number format + alignment) from the row above so the
# new cells reuse the existing date/time styles instead of minting new ones.
$ws.Range("B6:C6").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)

# New row of data (row 7) -- E7 ("1:00AM") is written before A7 so the two
# new shared-string entries are interned in the same order as the source.
$ws.Range("E7").Value = "1:00AM"
$ws.Range("A7").Value = "Forms Designing/Other Functions"
$ws.Range("B7").Value = 43761
$ws.Range("C7").Value = 43762
$ws.Range("D7").Value = 0.41666666666666669

# Widen column A slightly to fit the new text, and update the selection
$ws.Columns.Item(1).ColumnWidth = 30.5
$ws.Range("A8").Select()
